# Apply the "maicol" product row edit to Hoja1. Actividades
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1. Actividades")

# Row 10: clear the old "PRODUCTO" value ( factura) from J10
$ws.Range("J10").Value = $null

# Row 11: populate the new sample activity row about "maicol"
$ws.Range("A11").Value = "maicol es el mejor tipo que conozco"
$ws.Range("B11").Value = "maicol"
$ws.Range("C11").Value = "$$$"
$ws.Range("G11").Value = "X"
$ws.Range("I11").Value = "billetera"
$ws.Range("J11").Value = " Producto Maicol"
